$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 test case (was iAU_TC_ID_131 / negative-scenario text)
$ws.Range("A2").Value = " iAU_TC_ID_212"
$ws.Range("B2").Value = "@RegressionA Pre-Request Verify Elumina Login and Create Exam"
$ws.Range("C2").Value = "passed"

# New row 3
$ws.Range("A3").Value = " iAU_TC_ID_212"
$ws.Range("B3").Value = '@RegressionA Pre-Request "Validation of Delivery --> Add New Users"'
$ws.Range("C3").Value = "passed"

# New row 4
$ws.Range("A4").Value = "iAU_TC_ID_215"
$ws.Range("B4").Value = "@RegressionA Validation of Delivery--> Live Monitor all exam status  "
$ws.Range("C4").Value = "passed"
